$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("J:J").Delete()
$ws.Rows("3:4").Delete()

$row1 = @(1, 5, 2, 8, 6, 9, 3, 4, 7)
$row2 = @(635333980, 297083785, 141679615, 130143905, 115280370, 79968760, 70457945, 56205570, 27251910)

for ($i = 0; $i -lt 9; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $row1[$i]
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}

$ws.Range("A1:I2").NumberFormat = "0_);[Red](0)"
$ws.Rows("1:2").EntireRow.AutoFit()
